$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TPM-derived metrics for rows 2-10 (Slit1-Gpc1 LR pairs)
# as recomputed with the new TPM data.
$ws.Range("G2").Value = 0.04332866666666666
$ws.Range("H2").Value = 0.129986
$ws.Range("I2").Value = 0.2246397599897691
$ws.Range("J2").Value = 0.2246397599897691
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.09554499999999999
$ws.Range("N2").Value = 0.286635
$ws.Range("O2").Value = 0.01029975823317688
$ws.Range("P2").Value = 0.01029975823317688
$ws.Range("Q2").Value = 0.004139837456666666
$ws.Range("R2").Value = 0.03725853711
$ws.Range("S2").Value = 0.002313735217453503
$ws.Range("T2").Value = 0.002313735217453503
$ws.Range("G3").Value = 0.04332866666666666
$ws.Range("H3").Value = 0.129986
$ws.Range("I3").Value = 0.2246397599897691
$ws.Range("J3").Value = 0.2246397599897691
$ws.Range("O3").Value = 0.4011437372432085
$ws.Range("P3").Value = 0.4011437372432086
$ws.Range("Q3").Value = 0.1612338689268889
$ws.Range("R3").Value = 1.451104820342
$ws.Range("S3").Value = 0.09011283285571337
$ws.Range("T3").Value = 0.09011283285571339
$ws.Range("G4").Value = 0.04332866666666666
$ws.Range("H4").Value = 0.129986
$ws.Range("I4").Value = 0.2246397599897691
$ws.Range("J4").Value = 0.2246397599897691
$ws.Range("M4").Value = 5.459703999999999
$ws.Range("N4").Value = 16.379112
$ws.Range("O4").Value = 0.5885565045236145
$ws.Range("P4").Value = 0.5885565045236146
$ws.Range("Q4").Value = 0.2365616947146666
$ws.Range("R4").Value = 2.129055252432
$ws.Range("S4").Value = 0.1322131919166022
$ws.Range("T4").Value = 0.1322131919166023
$ws.Range("I5").Value = 0.5955530362469368
$ws.Range("J5").Value = 0.5955530362469369
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.09554499999999999
$ws.Range("N5").Value = 0.286635
$ws.Range("O5").Value = 0.01029975823317688
$ws.Range("P5").Value = 0.01029975823317688
$ws.Range("Q5").Value = 0.01097531784666667
$ws.Range("R5").Value = 0.09877786062
$ws.Range("S5").Value = 0.006134052288377877
$ws.Range("T5").Value = 0.006134052288377878
$ws.Range("I6").Value = 0.5955530362469368
$ws.Range("J6").Value = 0.5955530362469369
$ws.Range("O6").Value = 0.4011437372432085
$ws.Range("P6").Value = 0.4011437372432086
$ws.Range("S6").Value = 0.2389023706866363
$ws.Range("T6").Value = 0.2389023706866363
$ws.Range("I7").Value = 0.5955530362469368
$ws.Range("J7").Value = 0.5955530362469369
$ws.Range("M7").Value = 5.459703999999999
$ws.Range("N7").Value = 16.379112
$ws.Range("O7").Value = 0.5885565045236145
$ws.Range("P7").Value = 0.5885565045236146
$ws.Range("Q7").Value = 0.6271598382826666
$ws.Range("R7").Value = 5.644438544544
$ws.Range("S7").Value = 0.3505166132719226
$ws.Range("T7").Value = 0.3505166132719227
$ws.Range("G8").Value = 0.03468133333333333
$ws.Range("H8").Value = 0.104044
$ws.Range("I8").Value = 0.1798072037632941
$ws.Range("J8").Value = 0.1798072037632941
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.09554499999999999
$ws.Range("N8").Value = 0.286635
$ws.Range("O8").Value = 0.01029975823317688
$ws.Range("P8").Value = 0.01029975823317688
$ws.Range("Q8").Value = 0.003313627993333333
$ws.Range("R8").Value = 0.02982265194
$ws.Range("S8").Value = 0.001851970727345501
$ws.Range("T8").Value = 0.001851970727345501
$ws.Range("G9").Value = 0.03468133333333333
$ws.Range("H9").Value = 0.104044
$ws.Range("I9").Value = 0.1798072037632941
$ws.Range("J9").Value = 0.1798072037632941
$ws.Range("O9").Value = 0.4011437372432085
$ws.Range("P9").Value = 0.4011437372432086
$ws.Range("Q9").Value = 0.1290555648964444
$ws.Range("R9").Value = 1.161500084068
$ws.Range("S9").Value = 0.07212853370085889
$ws.Range("T9").Value = 0.07212853370085889
$ws.Range("G10").Value = 0.03468133333333333
$ws.Range("H10").Value = 0.104044
$ws.Range("I10").Value = 0.1798072037632941
$ws.Range("J10").Value = 0.1798072037632941
$ws.Range("M10").Value = 5.459703999999999
$ws.Range("N10").Value = 16.379112
$ws.Range("O10").Value = 0.5885565045236145
$ws.Range("P10").Value = 0.5885565045236146
$ws.Range("Q10").Value = 0.1893498143253333
$ws.Range("R10").Value = 1.704148328928
$ws.Range("S10").Value = 0.1058266993350897
$ws.Range("T10").Value = 0.1058266993350897